$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new time-log entry as row 18, right below the existing table.
$ws.Range("A18").Value = 46055
$ws.Range("B18").Value = 2

# Match the formatting used by the row above it (date style on A, plain
# number style on B) by copying the formats down from row 17.
$ws.Range("A17:B17").Copy() | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# The Time [h] total (D2) already contains =SUM(B:B); Excel recalculates it
# automatically as values change. Land the selection where Excel would
# leave it after typing the new row and pressing Enter.
$ws.Range("B19").Select()
